$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '27.865.87'
$cell.ClearFormats()
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.43%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.625.32'
$cell.ClearFormats()
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.99%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.19%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '211.04'
$cell.ClearFormats()
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.22%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.34%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.23%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '23.42'
$cell.ClearFormats()
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.90%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.19%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0612'
$cell.ClearFormats()
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.45%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.11%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.858.78'
$cell.ClearFormats()
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.86%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.621.25'
$cell.ClearFormats()
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.36%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.562'
$cell.ClearFormats()
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.44%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '65.44'
$cell.ClearFormats()
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.75%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '27.861.38'
$cell.ClearFormats()
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.43%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '230.22'
$cell.ClearFormats()
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.93%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.66'
$cell.ClearFormats()
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.48%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0721'
$cell.ClearFormats()
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.39%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.29%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.32'
$cell.ClearFormats()
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.95%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.14'
$cell.ClearFormats()
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -5.78%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.03'
$cell.ClearFormats()
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.56%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '154.63'
$cell.ClearFormats()

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.90'
$cell.ClearFormats()
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.01%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '15.53'
$cell.ClearFormats()
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.27%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.23%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.87%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0481'
$cell.ClearFormats()
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.73%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.41'
$cell.ClearFormats()
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.82%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.07'
$cell.ClearFormats()
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.04%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.398.86'
$cell.ClearFormats()
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.11%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.44%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +8.91%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.24%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0170'
$cell.ClearFormats()
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.04%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.37%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.862'
$cell.ClearFormats()
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.97%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.61%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.31%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.83'
$cell.ClearFormats()
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.50%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.23%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '65.79'
$cell.ClearFormats()
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.73%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.92%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.768.38'
$cell.ClearFormats()
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.84%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '88.04'
$cell.ClearFormats()
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.102'
$cell.ClearFormats()
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.63%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.38%  '
$cell.ClearFormats()

$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.47%  '
$cell.ClearFormats()
